$d = $word.ActiveDocument

# The phrase "третьего эксперта" occurs twice in this document (two
# near-duplicate paragraphs). The edit targets the SECOND occurrence
# (the one immediately followed by the "У вас не получилось узнать..."
# paragraph), so skip past the first match before searching again.

$first = $d.Content
$first.Find.Execute("третьего эксперта", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$first.Collapse(0)

$target = $d.Range($first.End, $d.Content.End)
$target.Find.Execute("третьего", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Replace "третьего" with "своего" while forcing Word to split the
# surrounding run into three runs (same formatting, just separate
# <w:r> elements) the way a real edit-in-place would: toggle a format
# property across the replacement so the run can't silently re-merge
# with its neighbours, then toggle it back to match the original
# (identical) formatting.
$target.Bold = 1
$target.Text = "своего"
$target.Bold = 0
